{"js": "// 1. Merge the \"Fridlysta arter\" intro paragraph and the two bullet\n//    list items into a single paragraph with combined text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nlet introPara = null;\nlet bulletIndices = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er\") !== -1) {\n    introPara = p;\n  } else if (introPara !== null && p.style === \"List Bullet\") {\n    bulletIndices.push(i);\n  }\n}\n\nintroPara.insertText(\n  \"F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er och v\u00e4xtplatser i den avverkningsanm\u00e4lda skogen: skogsfru (NT, \u00a78) och skogsr\u00f6r (\u00a77).\",\n  Word.InsertLocation.replace\n);\n\n// Delete the now-redundant bullet paragraphs (in reverse order so the\n// indices captured above stay valid while deleting).\nfor (let i = bulletIndices.length - 1; i >= 0; i--) {\n  paragraphs.items[bulletIndices[i]].delete();\n}\n\n// 2. Update the date stamped in the first-page header.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst firstPageHeader = sections.items[0].getHeader(\"FirstPage\");\nconst dateResults = firstPageHeader.search(\"2023-10-22\");\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2023-10-25\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Merge the \"Fridlysta arter\" intro paragraph and the two bullet\n#    list items below it into a single paragraph with combined text.\n$introIndex = 0\n$bulletIndices = @()\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er*\") {\n        $introIndex = $i\n    } elseif ($introIndex -ne 0 -and $p.Style.NameLocal -eq \"List Bullet\") {\n        $bulletIndices += $i\n    }\n    $i++\n}\n\n$d.Paragraphs($introIndex).Range.Text = \"F\u00f6ljande fridlysta arter har sina livsmilj\u00f6er och v\u00e4xtplatser i den avverkningsanm\u00e4lda skogen: skogsfru (NT, \u00a78) och skogsr\u00f6r (\u00a77).\"\n\n# Delete the bullet paragraphs from the highest index down so earlier\n# indices remain valid.\n$sortedBullets = $bulletIndices | Sort-Object -Descending\nforeach ($idx in $sortedBullets) {\n    $d.Paragraphs($idx).Range.Delete()\n}\n\n# 2. Update the date stamped in the first-page header.\n$sec = $d.Sections(1)\n$firstPageHeader = $sec.Headers(2)\n$firstPageHeader.Range.Find.Execute(\"2023-10-22\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-10-25\", 2)\n"}
